# Apply price/coin-listing updates from the "Updated symbol list" GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.20"
$ws.Range("D3").Value = "'23.09"
$ws.Range("D4").Value = "'5.412"
$ws.Range("D5").Value = "'0.05983"
$ws.Range("D6").Value = "'3.462"
$ws.Range("D8").Value = "'0.8139"
$ws.Range("D9").Value = "'0.9142"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1408"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07417"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03234"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03089"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09362"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.850"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001554"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04668"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005936"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006087"
$ws.Range("D20").Value = "'0.005005"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("D21").Value = "'0.0009878"
$ws.Range("D22").Value = "'0.00007798"
$ws.Range("D23").Value = "'3.614"
$ws.Range("D24").Value = "'2.130"
$ws.Range("D26").Value = "'0.1302"
$ws.Range("D27").Value = "'0.0002898"
$ws.Range("D40").Value = "'0.03928"
$ws.Range("D41").Value = "'0.006208"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("D43").Value = "'0.002619"
$ws.Range("D44").Value = "'0.006439"
$ws.Range("D45").Value = "'0.00005242"
$ws.Range("D48").Value = "'0.8295"
$ws.Range("D49").Value = "'0.002270"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D51").Value = "'0.0001999"
